# Auto-generated: apply scraped price/profit refresh to the Leve profit tables.
# Each worksheet corresponds to a crafting class; columns H-N are the
# market-driven columns (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 56.166668
$ws.Range("I11").Value = 56.166668
$ws.Range("K11").Value = 56.166668
$ws.Range("M11").Value = 83.833332
$ws.Range("H18").Value = 1027.8125
$ws.Range("I18").Value = 997
$ws.Range("K18").Value = 997
$ws.Range("M18").Value = -713
$ws.Range("H43").Value = 1449.5
$ws.Range("I43").Value = 1450
$ws.Range("J43").Value = 1449
$ws.Range("K43").Value = 1450
$ws.Range("L43").Value = 1449
$ws.Range("M43").Value = -1381
$ws.Range("N43").Value = -1587
$ws.Range("H80").Value = 1031.2222
$ws.Range("J80").Value = 1456.2
$ws.Range("L80").Value = 4368.6
$ws.Range("N80").Value = -6364.6
$ws.Range("H83").Value = 1031.2222
$ws.Range("J83").Value = 1456.2
$ws.Range("L83").Value = 13105.8
$ws.Range("N83").Value = -23089.8
$ws.Range("H86").Value = 488646.25
$ws.Range("J86").Value = 488646.25
$ws.Range("L86").Value = 488646.25
$ws.Range("N86").Value = -490892.25
$ws.Range("H89").Value = 488646.25
$ws.Range("J89").Value = 488646.25
$ws.Range("L89").Value = 2443231.25
$ws.Range("N89").Value = -2454463.25
$ws.Range("H92").Value = 1518
$ws.Range("I92").Value = 1499
$ws.Range("J92").Value = 1527.5
$ws.Range("K92").Value = 1499
$ws.Range("L92").Value = 1527.5
$ws.Range("M92").Value = -251
$ws.Range("N92").Value = -4023.5
$ws.Range("H111").Value = 5939
$ws.Range("I111").Value = 6071.067
$ws.Range("J111").Value = 5542.8
$ws.Range("K111").Value = 18213.201
$ws.Range("L111").Value = 16628.4
$ws.Range("M111").Value = -15146.201
$ws.Range("N111").Value = -22762.4
$ws.Range("H112").Value = 4347.222
$ws.Range("J112").Value = 4440.625
$ws.Range("L112").Value = 13321.875
$ws.Range("N112").Value = -15537.875
$ws.Range("H123").Value = 69999.5
$ws.Range("J123").Value = 69999.5
$ws.Range("L123").Value = 69999.5
$ws.Range("N123").Value = -79799.5
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H137").Value = 3599.3
$ws.Range("I137").Value = 1427.8572
$ws.Range("K137").Value = 4283.571599999999
$ws.Range("M137").Value = -1733.571599999999
$ws.Range("H138").Value = 2446.15
$ws.Range("I138").Value = 1786.3572
$ws.Range("J138").Value = 3985.6667
$ws.Range("K138").Value = 5359.071599999999
$ws.Range("L138").Value = 11957.0001
$ws.Range("M138").Value = -219.0715999999993
$ws.Range("N138").Value = -22237.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1418.125
$ws.Range("I61").Value = 1090.8334
$ws.Range("K61").Value = 1090.8334
$ws.Range("M61").Value = -878.8334
$ws.Range("H74").Value = 1684.25
$ws.Range("I74").Value = 1699
$ws.Range("K74").Value = 1699
$ws.Range("M74").Value = -825
$ws.Range("H77").Value = 1684.25
$ws.Range("I77").Value = 1699
$ws.Range("K77").Value = 8495
$ws.Range("M77").Value = -4127
$ws.Range("H110").Value = 1204.1111
$ws.Range("I110").Value = 1309.1333
$ws.Range("K110").Value = 1309.1333
$ws.Range("M110").Value = 735.8667
$ws.Range("H132").Value = 931.5714
$ws.Range("J132").Value = 966.3333
$ws.Range("L132").Value = 2898.9999
$ws.Range("N132").Value = -7958.9999
$ws.Range("H136").Value = 1418.125
$ws.Range("I136").Value = 1090.8334
$ws.Range("K136").Value = 3272.5002
$ws.Range("M136").Value = -722.5001999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 206.22223
$ws.Range("J22").Value = 200.75
$ws.Range("L22").Value = 200.75
$ws.Range("N22").Value = -546.75
$ws.Range("H86").Value = 1448.875
$ws.Range("I86").Value = 1298.5555
$ws.Range("J86").Value = 1899.8334
$ws.Range("K86").Value = 1298.5555
$ws.Range("L86").Value = 1899.8334
$ws.Range("M86").Value = -175.5554999999999
$ws.Range("N86").Value = -4145.8334
$ws.Range("H89").Value = 1448.875
$ws.Range("I89").Value = 1298.5555
$ws.Range("J89").Value = 1899.8334
$ws.Range("K89").Value = 6492.7775
$ws.Range("L89").Value = 9499.166999999999
$ws.Range("M89").Value = -876.7775000000001
$ws.Range("N89").Value = -20731.167
$ws.Range("H99").Value = 1713.4286
$ws.Range("I99").Value = 1665.6666
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1665.6666
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -167.6666
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 1921.6666
$ws.Range("I134").Value = 1921.6666
$ws.Range("K134").Value = 5764.9998
$ws.Range("M134").Value = -3229.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1507
$ws.Range("I31").Value = 1507
$ws.Range("K31").Value = 1507
$ws.Range("M31").Value = -1212
$ws.Range("H34").Value = 1507
$ws.Range("I34").Value = 1507
$ws.Range("K34").Value = 1507
$ws.Range("M34").Value = -1305
$ws.Range("H107").Value = 846.5625
$ws.Range("I107").Value = 548.3
$ws.Range("K107").Value = 548.3
$ws.Range("M107").Value = 1371.7
$ws.Range("H132").Value = 3097.6667
$ws.Range("I132").Value = 2917.8
$ws.Range("K132").Value = 8753.400000000001
$ws.Range("M132").Value = -6223.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 30000140
$ws.Range("I23").Value = 50000090
$ws.Range("J23").Value = 218.5
$ws.Range("K23").Value = 150000270
$ws.Range("L23").Value = 655.5
$ws.Range("M23").Value = -150000035
$ws.Range("N23").Value = -1125.5
$ws.Range("H70").Value = 11527.333
$ws.Range("I70").Value = 1374.5
$ws.Range("J70").Value = 14428.143
$ws.Range("K70").Value = 4123.5
$ws.Range("L70").Value = 43284.429
$ws.Range("M70").Value = -3808.5
$ws.Range("N70").Value = -43914.429
$ws.Range("H73").Value = 11527.333
$ws.Range("I73").Value = 1374.5
$ws.Range("J73").Value = 14428.143
$ws.Range("K73").Value = 4123.5
$ws.Range("L73").Value = 43284.429
$ws.Range("M73").Value = -3031.5
$ws.Range("N73").Value = -45468.429
$ws.Range("H97").Value = 1574.5
$ws.Range("I97").Value = 2599
$ws.Range("K97").Value = 7797
$ws.Range("M97").Value = -7301

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 51500
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 80000
$ws.Range("N57").Value = -81640
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
$ws.Range("H122").Value = 4607.3076
$ws.Range("I122").Value = 4099.75
$ws.Range("J122").Value = 4832.8887
$ws.Range("K122").Value = 12299.25
$ws.Range("L122").Value = 14498.6661
$ws.Range("M122").Value = -9849.25
$ws.Range("N122").Value = -19398.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1961
$ws.Range("I22").Value = 1837.6316
$ws.Range("K22").Value = 1837.6316
$ws.Range("M22").Value = -1542.6316
$ws.Range("H27").Value = 1961
$ws.Range("I27").Value = 1837.6316
$ws.Range("K27").Value = 1837.6316
$ws.Range("M27").Value = -1730.6316
$ws.Range("H46").Value = 1874.6428
$ws.Range("I46").Value = 1343.8889
$ws.Range("J46").Value = 2830
$ws.Range("K46").Value = 1343.8889
$ws.Range("L46").Value = 2830
$ws.Range("M46").Value = -1155.8889
$ws.Range("N46").Value = -3206
$ws.Range("H82").Value = 1726.8182
$ws.Range("I82").Value = 1139.8
$ws.Range("J82").Value = 2216
$ws.Range("K82").Value = 1139.8
$ws.Range("L82").Value = 2216
$ws.Range("M82").Value = -778.8
$ws.Range("N82").Value = -2938
$ws.Range("H85").Value = 1726.8182
$ws.Range("I85").Value = 1139.8
$ws.Range("J85").Value = 2216
$ws.Range("K85").Value = 1139.8
$ws.Range("L85").Value = 2216
$ws.Range("M85").Value = 108.2
$ws.Range("N85").Value = -4712
$ws.Range("H93").Value = 2000.3334
$ws.Range("I93").Value = 2000.3334
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2000.3334
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -752.3334
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 4441.3335
$ws.Range("I100").Value = 4769.8
$ws.Range("J100").Value = 2799
$ws.Range("K100").Value = 4769.8
$ws.Range("L100").Value = 2799
$ws.Range("M100").Value = -4228.8
$ws.Range("N100").Value = -3881
$ws.Range("H132").Value = 2032.3529
$ws.Range("I132").Value = 1638.0834
$ws.Range("K132").Value = 4914.2502
$ws.Range("M132").Value = -2384.2502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8024.5557
$ws.Range("I122").Value = 6890.636
$ws.Range("K122").Value = 20671.908
$ws.Range("M122").Value = -18221.908
$ws.Range("H132").Value = 2072.875
$ws.Range("I132").Value = 2072.875
$ws.Range("K132").Value = 6218.625
$ws.Range("M132").Value = -3688.625
$ws.Range("H136").Value = 887.8461
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = -150
